# Stückliste Link update
# - Adds hyperlinks to the distributor "Link" column (L2:L4)
# - Removes the blank spacer row above the "Summe" total row and adds a
#   thin top border / rule above the total row
# - Updates the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the empty spacer row (old row 10); this shifts the "Summe"
#     total row from 11 up to 10, while the SUM(J2:J9) formula keeps
#     pointing at the same cells.
$ws.Rows(10).Delete()

# --- Draw a thin rule above the new "Summe" row (A10:K10)
$sumBorderRange = $ws.Range("A10:K10")
$topBorder = $sumBorderRange.Borders.Item(8)   # xlEdgeTop
$topBorder.LineStyle = 1                        # xlContinuous
$topBorder.Weight = 2                           # xlThin

# --- Turn the distributor links into real hyperlinks
$linkRsOnline = "https://at.rs-online.com/web/p/operationsverstarker/7591534/?relevancy-data=636F3D3126696E3D4931384E53656172636847656E65726963266C753D6465266D6D3D6D61746368616C6C7061727469616C26706D3D5E5B5C707B4C7D5C707B4E647D2D2C2F255C2E5D2B2426706F3D31333326736E3D592673723D2673743D4B4559574F52445F53494E474C455F414C5048415F4E554D455249432673633D592677633D4E4F4E45267573743D414441343839312D3241525A267374613D414441343839312D3241525A26&searchHistory=%7B%22enabled%22%3Atrue%7D"
$linkFarnellAdapter = "https://at.farnell.com/aries/lcqt-soic8-8/ic-adapter-8-soic-dip-2-54mm/dp/2476033?st=SOIC%20Adapter"
$linkFarnellLm2776 = "https://at.farnell.com/texas-instruments/lm2776dbvr/dc-dc-ladungspumpe-invertierend/dp/2817376?st=LM2776DBVR"

$cellL2 = $ws.Range("L2")
$ws.Hyperlinks.Add($cellL2, $linkRsOnline, "", "", $cellL2.Value2)

$ws.Hyperlinks.Add($ws.Range("L3"), $linkFarnellAdapter)

$ws.Hyperlinks.Add($ws.Range("L4"), $linkFarnellLm2776)

# --- Restore the active selection
$ws.Range("E24").Select()
